# "switched left & right again"
# Swap the "left"/"right" labels between the two mic-plane blocks on
# worksheet "Tabelle1":
#   - rows 2-18  (plane B) currently say "left"  in columns C & D -> "right"
#   - rows 36-52 (plane D) currently say "right" in columns C & D -> "left"
# (The last row of each block, 18 & 52, only has the C column changed
#  because its D column already holds "top".)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = "right"   # column C
    $ws.Cells.Item($r, 4).Value = "right"   # column D
}
# Row 18's column D already contains "top" and must stay that way.
$ws.Cells.Item(18, 4).Value = "top"

for ($r = 36; $r -le 52; $r++) {
    $ws.Cells.Item($r, 3).Value = "left"    # column C
    $ws.Cells.Item($r, 4).Value = "left"    # column D
}
# Row 52's column D already contains "top" and must stay that way.
$ws.Cells.Item(52, 4).Value = "top"

# Update the view state to match the saved selection/scroll position.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$ws.Range("E43").Select()
